$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62,8).Value = 2654.8948
$ws.Cells.Item(62,9).Value = 1949.75
$ws.Cells.Item(62,10).Value = 3863.7144
$ws.Cells.Item(62,11).Value = 1949.75
$ws.Cells.Item(62,12).Value = 3863.7144
$ws.Cells.Item(62,13).Value = -1325.75
$ws.Cells.Item(62,14).Value = -5111.7144

$ws.Cells.Item(65,8).Value = 2654.8948
$ws.Cells.Item(65,9).Value = 1949.75
$ws.Cells.Item(65,10).Value = 3863.7144
$ws.Cells.Item(65,11).Value = 9748.75
$ws.Cells.Item(65,12).Value = 19318.572
$ws.Cells.Item(65,13).Value = -6628.75
$ws.Cells.Item(65,14).Value = -25558.572

$ws.Cells.Item(129,8).Value = 250793
$ws.Cells.Item(129,10).Value = 278617.22
$ws.Cells.Item(129,12).Value = 835851.6599999999
$ws.Cells.Item(129,14).Value = -845851.6599999999

$ws.Cells.Item(138,8).Value = 1358.61
$ws.Cells.Item(138,9).Value = 557.0227
$ws.Cells.Item(138,10).Value = 1988.4286
$ws.Cells.Item(138,11).Value = 1671.0681
$ws.Cells.Item(138,12).Value = 5965.2858
$ws.Cells.Item(138,13).Value = 3468.9319
$ws.Cells.Item(138,14).Value = -16245.2858

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2,8).Value = 1592.7391
$ws.Cells.Item(2,9).Value = 1369.4736
$ws.Cells.Item(2,10).Value = 2653.25
$ws.Cells.Item(2,11).Value = 1369.4736
$ws.Cells.Item(2,12).Value = 2653.25
$ws.Cells.Item(2,13).Value = -1256.4736
$ws.Cells.Item(2,14).Value = -2879.25

$ws.Cells.Item(32,8).Value = 34299.152
$ws.Cells.Item(32,9).Value = 34276.812
$ws.Cells.Item(32,10).Value = 35014
$ws.Cells.Item(32,11).Value = 34276.812
$ws.Cells.Item(32,12).Value = 35014
$ws.Cells.Item(32,13).Value = -33989.812
$ws.Cells.Item(32,14).Value = -35588

$ws.Cells.Item(116,8).Value = 1592.7391
$ws.Cells.Item(116,9).Value = 1369.4736
$ws.Cells.Item(116,10).Value = 2653.25
$ws.Cells.Item(116,11).Value = 1369.4736
$ws.Cells.Item(116,12).Value = 2653.25
$ws.Cells.Item(116,13).Value = 924.5264
$ws.Cells.Item(116,14).Value = -7241.25

$ws.Cells.Item(134,8).Value = 52604.75
$ws.Cells.Item(134,10).Value = 52604.75
$ws.Cells.Item(134,12).Value = 52604.75
$ws.Cells.Item(134,14).Value = -62744.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3,8).Value = 1592.7391
$ws.Cells.Item(3,9).Value = 1369.4736
$ws.Cells.Item(3,10).Value = 2653.25
$ws.Cells.Item(3,11).Value = 1369.4736
$ws.Cells.Item(3,12).Value = 2653.25
$ws.Cells.Item(3,13).Value = -1255.4736
$ws.Cells.Item(3,14).Value = -2881.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58,8).Value = 15375.571
$ws.Cells.Item(58,9).Value = 1184.9565
$ws.Cells.Item(58,10).Value = 42574.25
$ws.Cells.Item(58,11).Value = 1184.9565
$ws.Cells.Item(58,12).Value = 42574.25
$ws.Cells.Item(58,13).Value = -981.9565
$ws.Cells.Item(58,14).Value = -42980.25

$ws.Cells.Item(86,8).Value = 5959882.5
$ws.Cells.Item(86,9).Value = 2580.8823
$ws.Cells.Item(86,11).Value = 2580.8823
$ws.Cells.Item(86,13).Value = -1457.8823

$ws.Cells.Item(89,8).Value = 5959882.5
$ws.Cells.Item(89,9).Value = 2580.8823
$ws.Cells.Item(89,11).Value = 12904.4115
$ws.Cells.Item(89,13).Value = -7288.411500000002

$ws.Cells.Item(99,8).Value = 15155220
$ws.Cells.Item(99,9).Value = 3090.7368
$ws.Cells.Item(99,10).Value = 35718824
$ws.Cells.Item(99,11).Value = 3090.7368
$ws.Cells.Item(99,12).Value = 35718824
$ws.Cells.Item(99,13).Value = -1592.7368
$ws.Cells.Item(99,14).Value = -35721820

$ws.Cells.Item(126,8).Value = 15155220
$ws.Cells.Item(126,9).Value = 3090.7368
$ws.Cells.Item(126,10).Value = 35718824
$ws.Cells.Item(126,11).Value = 9272.2104
$ws.Cells.Item(126,12).Value = 107156472
$ws.Cells.Item(126,13).Value = -6802.2104
$ws.Cells.Item(126,14).Value = -107161412

$ws.Cells.Item(134,8).Value = 1139.9434
$ws.Cells.Item(134,9).Value = 832.8077
$ws.Cells.Item(134,11).Value = 2498.4231
$ws.Cells.Item(134,13).Value = 36.57690000000002

$ws.Cells.Item(136,8).Value = 15375.571
$ws.Cells.Item(136,9).Value = 1184.9565
$ws.Cells.Item(136,10).Value = 42574.25
$ws.Cells.Item(136,11).Value = 3554.8695
$ws.Cells.Item(136,12).Value = 127722.75
$ws.Cells.Item(136,13).Value = -1004.8695
$ws.Cells.Item(136,14).Value = -132822.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5,8).Value = 1197.625
$ws.Cells.Item(5,9).Value = 1062
$ws.Cells.Item(5,11).Value = 3186
$ws.Cells.Item(5,13).Value = -3074

$ws.Cells.Item(107,8).Value = 7948.5386
$ws.Cells.Item(107,9).Value = 16882.166
$ws.Cells.Item(107,10).Value = 291.14285
$ws.Cells.Item(107,11).Value = 50646.49800000001
$ws.Cells.Item(107,12).Value = 873.4285500000001
$ws.Cells.Item(107,13).Value = -48726.49800000001
$ws.Cells.Item(107,14).Value = -4713.428550000001

$ws.Cells.Item(131,8).Value = 773.14
$ws.Cells.Item(131,9).Value = 633.3333
$ws.Cells.Item(131,10).Value = 777.4639
$ws.Cells.Item(131,11).Value = 1899.9999
$ws.Cells.Item(131,12).Value = 2332.3917
$ws.Cells.Item(131,13).Value = 3140.0001
$ws.Cells.Item(131,14).Value = -12412.3917

$ws.Cells.Item(132,8).Value = 1456.8
$ws.Cells.Item(132,9).Value = 1413
$ws.Cells.Item(132,11).Value = 12717
$ws.Cells.Item(132,13).Value = -10187

$ws.Cells.Item(135,8).Value = 1197.625
$ws.Cells.Item(135,9).Value = 1062
$ws.Cells.Item(135,11).Value = 9558
$ws.Cells.Item(135,13).Value = -7023

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(44,8).Value = 19999.5
$ws.Cells.Item(44,10).Value = 19999.5
$ws.Cells.Item(44,12).Value = 19999.5
$ws.Cells.Item(44,14).Value = -21191.5

$ws.Cells.Item(52,8).Value = 16672500
$ws.Cells.Item(52,10).Value = 16672500
$ws.Cells.Item(52,12).Value = 16672500
$ws.Cells.Item(52,14).Value = -16673018

$ws.Cells.Item(132,8).Value = 53925.066
$ws.Cells.Item(132,9).Value = 54394.35
$ws.Cells.Item(132,10).Value = 52986.5
$ws.Cells.Item(132,11).Value = 163183.05
$ws.Cells.Item(132,12).Value = 158959.5
$ws.Cells.Item(132,13).Value = -160653.05
$ws.Cells.Item(132,14).Value = -164019.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(94,8).Value = 39665
$ws.Cells.Item(94,10).Value = 39665
$ws.Cells.Item(94,12).Value = 39665
$ws.Cells.Item(94,14).Value = -41017

$ws.Cells.Item(136,8).Value = 46509.184
$ws.Cells.Item(136,9).Value = 46509.184
$ws.Cells.Item(136,10).Value = 0
$ws.Cells.Item(136,11).Value = 139527.552
$ws.Cells.Item(136,12).Value = 0
$ws.Cells.Item(136,13).Value = -136977.552
$ws.Cells.Item(136,14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113,8).Value = 2458102
$ws.Cells.Item(113,9).Value = 1833.3334
$ws.Cells.Item(113,10).Value = 5405624
$ws.Cells.Item(113,11).Value = 5500.0002
$ws.Cells.Item(113,12).Value = 16216872
$ws.Cells.Item(113,13).Value = -3330.0002
$ws.Cells.Item(113,14).Value = -16221212

$ws.Cells.Item(132,8).Value = 1166.7587
$ws.Cells.Item(132,9).Value = 921.381
$ws.Cells.Item(132,10).Value = 1810.875
$ws.Cells.Item(132,11).Value = 2764.143
$ws.Cells.Item(132,12).Value = 5432.625
$ws.Cells.Item(132,13).Value = -234.143
$ws.Cells.Item(132,14).Value = -10492.625

$ws.Cells.Item(136,8).Value = 26317592
$ws.Cells.Item(136,9).Value = 50001748
$ws.Cells.Item(136,10).Value = 1861.5555
$ws.Cells.Item(136,11).Value = 150005244
$ws.Cells.Item(136,12).Value = 5584.666499999999
$ws.Cells.Item(136,13).Value = -150002694
$ws.Cells.Item(136,14).Value = -10684.6665
